$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'296.94"
$ws.Cells.Item(2, 5).Value = "'-1.10%"
$ws.Cells.Item(3, 4).Value = "'31.41"
$ws.Cells.Item(3, 5).Value = "'0.08%"
$ws.Cells.Item(4, 4).Value = "'5.075"
$ws.Cells.Item(4, 5).Value = "'-1.27%"
$ws.Cells.Item(5, 4).Value = "'0.07990"
$ws.Cells.Item(5, 5).Value = "'8.82%"
$ws.Cells.Item(6, 4).Value = "'2.482"
$ws.Cells.Item(6, 5).Value = "'36.94%"
$ws.Cells.Item(7, 4).Value = "'7.775"
$ws.Cells.Item(7, 5).Value = "'-0.12%"
$ws.Cells.Item(8, 4).Value = "'0.9254"
$ws.Cells.Item(8, 5).Value = "'-0.03%"
$ws.Cells.Item(9, 4).Value = "'0.1746"
$ws.Cells.Item(9, 5).Value = "'3.79%"
$ws.Cells.Item(10, 4).Value = "'0.07376"
$ws.Cells.Item(10, 5).Value = "'3.08%"
$ws.Cells.Item(11, 5).Value = "'10.52%"
$ws.Cells.Item(12, 4).Value = "'0.03030"
$ws.Cells.Item(12, 5).Value = "'-0.46%"
$ws.Cells.Item(13, 4).Value = "'0.1001"
$ws.Cells.Item(13, 5).Value = "'0.77%"
$ws.Cells.Item(14, 4).Value = "'0.001492"
$ws.Cells.Item(14, 5).Value = "'0.25%"
$ws.Cells.Item(15, 4).Value = "'0.005963"
$ws.Cells.Item(15, 5).Value = "'-3.14%"
$ws.Cells.Item(16, 4).Value = "'3.519"
$ws.Cells.Item(16, 5).Value = "'1.75%"
$ws.Cells.Item(17, 5).Value = "'1.67%"
$ws.Cells.Item(18, 4).Value = "'2.243"
$ws.Cells.Item(18, 5).Value = "'0.90%"
$ws.Cells.Item(19, 5).Value = "'0.86%"
$ws.Cells.Item(20, 4).Value = "'0.1336"
$ws.Cells.Item(20, 5).Value = "'0.54%"
$ws.Cells.Item(21, 4).Value = "'4.291"
$ws.Cells.Item(21, 5).Value = "'-5.69%"
$ws.Cells.Item(22, 4).Value = "'0.1617"
$ws.Cells.Item(22, 5).Value = "'2.27%"
$ws.Cells.Item(23, 4).Value = "'0.04606"
$ws.Cells.Item(23, 5).Value = "'-0.97%"
$ws.Cells.Item(24, 4).Value = "'0.001246"
$ws.Cells.Item(24, 5).Value = "'2.70%"
$ws.Cells.Item(25, 4).Value = "'0.004421"
$ws.Cells.Item(25, 5).Value = "'-6.82%"
$ws.Cells.Item(26, 4).Value = "'0.0001201"
$ws.Cells.Item(26, 5).Value = "'-7.41%"
$ws.Cells.Item(27, 4).Value = "'0.0003427"
$ws.Cells.Item(27, 5).Value = "'83.09%"
$ws.Cells.Item(39, 4).Value = "'0.01770"
$ws.Cells.Item(39, 5).Value = "'2.64%"
$ws.Cells.Item(40, 4).Value = "'0.04468"
$ws.Cells.Item(40, 5).Value = "'-0.47%"
$ws.Cells.Item(41, 4).Value = "'0.006876"
$ws.Cells.Item(41, 5).Value = "'-3.09%"
$ws.Cells.Item(42, 5).Value = "'0.47%"
$ws.Cells.Item(43, 4).Value = "'0.002147"
$ws.Cells.Item(43, 5).Value = "'-3.06%"
$ws.Cells.Item(44, 4).Value = "'0.009790"
$ws.Cells.Item(44, 5).Value = "'-7.31%"
$ws.Cells.Item(45, 4).Value = "'0.00006564"
$ws.Cells.Item(45, 5).Value = "'5.05%"
$ws.Cells.Item(46, 4).Value = "'0.00000000751"
$ws.Cells.Item(46, 5).Value = "'0.07%"
$ws.Cells.Item(48, 5).Value = "'-55.54%"
$ws.Cells.Item(49, 4).Value = "'0.00002102"
$ws.Cells.Item(49, 5).Value = "'0.07%"
$ws.Cells.Item(50, 4).Value = "'0.0002001"
$ws.Cells.Item(50, 5).Value = "'0.14%"
